# Updates cryptos list with the latest scraped price/volume data.
# Note: several "Price" column values look numeric (e.g. "245.30", "0.661")
# but must stay as plain text, exactly like the source data (trailing
# zeros, fixed decimal places, etc.). A leading apostrophe forces Excel to
# store such values as text instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.357.57'
$ws.Range('E2').Value = '  -1.55%  '
$ws.Range('D3').Value = '2.051.66'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''245.30'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').Value = '''0.661'
$ws.Range('E6').Value = '  +1.09%  '
$ws.Range('D7').Value = '''56.80'
$ws.Range('E7').Value = '  +1.47%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '''63.08'
$ws.Range('E9').Value = '  +5.15%  '
$ws.Range('D10').Value = '''0.369'
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('D11').Value = '''0.0749'
$ws.Range('E11').Value = '  -3.24%  '
$ws.Range('E12').Value = '  -3.58%  '
$ws.Range('D13').Value = '''0.926'
$ws.Range('E13').Value = '  +3.52%  '
$ws.Range('D14').Value = '''14.41'
$ws.Range('E14').Value = '  -5.02%  '
$ws.Range('D15').Value = '2.348.81'
$ws.Range('E15').Value = '  -2.39%  '
$ws.Range('D16').Value = '''5.41'
$ws.Range('E16').Value = '  -2.79%  '
$ws.Range('D17').Value = '2.037.21'
$ws.Range('E17').Value = '  -3.29%  '
$ws.Range('D18').Value = '''17.80'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('D19').Value = '36.350.82'
$ws.Range('E19').Value = '  -1.52%  '
$ws.Range('D20').Value = '''71.61'
$ws.Range('E20').Value = '  -2.60%  '
$ws.Range('D21').Value = '0.0₃0858'
$ws.Range('E21').Value = '  -3.42%  '
$ws.Range('D22').Value = '''236.35'
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('D23').Value = '''5.23'
$ws.Range('E23').Value = '  -5.45%  '
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').Value = '''2.36'
$ws.Range('E25').Value = '  -2.71%  '
$ws.Range('D26').Value = '''2.27'
$ws.Range('E26').Value = '  +4.04%  '
$ws.Range('D27').Value = '''9.39'
$ws.Range('E27').Value = '  -5.34%  '
$ws.Range('D28').Value = '''164.67'
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('D29').Value = '''20.00'
$ws.Range('E29').Value = '  -4.00%  '
$ws.Range('E30').Value = '  -1.81%  '
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('D32').Value = '''5.00'
$ws.Range('E32').Value = '  -7.69%  '
$ws.Range('D33').Value = '''0.0600'
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('E34').Value = '  -6.28%  '
$ws.Range('D35').Value = '''0.0884'
$ws.Range('E35').Value = '  +4.51%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('D37').Value = '''1.84'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('E38').Value = '  -7.96%  '
$ws.Range('D39').Value = '''5.07'
$ws.Range('E39').Value = '  +3.29%  '
$ws.Range('E40').Value = '  -5.69%  '
$ws.Range('E41').Value = '  +0.90%  '
$ws.Range('E42').Value = '  -2.66%  '
$ws.Range('E43').Value = '  -5.64%  '
$ws.Range('E44').Value = '  -3.79%  '
$ws.Range('E45').Value = '  -4.83%  '
$ws.Range('D46').Value = '1.407.30'
$ws.Range('E46').Value = '  +3.53%  '
$ws.Range('D47').Value = '''15.95'
$ws.Range('E47').Value = '  -1.88%  '
$ws.Range('D48').Value = '''7.48'
$ws.Range('E48').Value = '  +6.15%  '
$ws.Range('D49').Value = '''2.96'
$ws.Range('E49').Value = '  +1.77%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '''46.30'
$ws.Range('E50').Value = '  +0.99%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '''2.27'
$ws.Range('E51').Value = '  -8.15%  '
